$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 29
$ws.Range("H29").Value = 309.2
$ws.Range("I29").Value = 260.5
$ws.Range("J29").Value = 504
$ws.Range("K29").Value = 781.5
$ws.Range("L29").Value = 1512
$ws.Range("M29").Value = -500.5
$ws.Range("N29").Value = -2074
# Row 69
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()
# Row 70
$ws.Range("H70").Value = 1833.3334
$ws.Range("J70").Value = 1833.3334
$ws.Range("L70").Value = 5500.0002
$ws.Range("N70").Value = -6040.0002
# Row 72
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()
# Row 73
$ws.Range("H73").Value = 1833.3334
$ws.Range("J73").Value = 1833.3334
$ws.Range("L73").Value = 5500.0002
$ws.Range("N73").Value = -7372.0002
# Row 135
$ws.Range("H135").Value = 1677
$ws.Range("J135").Value = 2181
$ws.Range("L135").Value = 19629
$ws.Range("N135").Value = -24699
# Row 138
$ws.Range("H138").Value = 2253.2917
$ws.Range("I138").Value = 760
$ws.Range("J138").Value = 2999.9375
$ws.Range("K138").Value = 2280
$ws.Range("L138").Value = 8999.8125
$ws.Range("M138").Value = 2860
$ws.Range("N138").Value = -19279.8125
# Row 141
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 10
$ws.Range("H10").Value = 300
$ws.Range("I10").Value = 300
$ws.Range("K10").Value = 300
$ws.Range("M10").Value = -130
# Row 24
$ws.Range("H24").Value = 50000
$ws.Range("J24").Value = 50000
$ws.Range("L24").Value = 50000
$ws.Range("N24").Value = -50748
# Row 32
$ws.Range("H32").Value = 3821.5264
$ws.Range("I32").Value = 4207.3335
$ws.Range("K32").Value = 4207.3335
$ws.Range("M32").Value = -3920.3335
# Row 100
$ws.Range("H100").Value = 50000
$ws.Range("J100").Value = 50000
$ws.Range("L100").Value = 50000
$ws.Range("N100").Value = -52164
# Row 132
$ws.Range("H132").Value = 1204
$ws.Range("I132").Value = 1204
$ws.Range("K132").Value = 3612
$ws.Range("M132").Value = -1082

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 102
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").ClearContents()
# Row 105
$ws.Range("H105").Value = 9899.5
$ws.Range("I105").Value = 9899
$ws.Range("K105").Value = 9899
$ws.Range("M105").Value = -8152

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 82
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
# Row 85
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
# Row 132
$ws.Range("H132").Value = 7765.5835
$ws.Range("I132").Value = 3946.75
$ws.Range("K132").Value = 11840.25
$ws.Range("M132").Value = -9310.25

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 17.777779
$ws.Range("I2").Value = 18.5
$ws.Range("J2").Value = 16.333334
$ws.Range("K2").Value = 111
$ws.Range("L2").Value = 98.000004
$ws.Range("M2").Value = 2
$ws.Range("N2").Value = -324.000004
# Row 10
$ws.Range("H10").Value = 378.25
$ws.Range("I10").Value = 4.3333335
$ws.Range("J10").Value = 1500
$ws.Range("K10").Value = 13.0000005
$ws.Range("L10").Value = 4500
$ws.Range("M10").Value = 125.9999995
$ws.Range("N10").Value = -4778
# Row 94
$ws.Range("H94").Value = 1000
$ws.Range("J94").Value = 1000
$ws.Range("L94").Value = 3000
$ws.Range("N94").Value = -4352
# Row 129
$ws.Range("H129").Value = 14998.5
$ws.Range("J129").Value = 14998.5
$ws.Range("L129").Value = 44995.5
$ws.Range("N129").Value = -54995.5
# Row 131
$ws.Range("H131").Value = 1483.3334
$ws.Range("I131").Value = 966.6667
$ws.Range("K131").Value = 2900.0001
$ws.Range("M131").Value = 2139.9999
# Row 137
$ws.Range("H137").Value = 18000
$ws.Range("J137").Value = 18000
$ws.Range("L137").Value = 54000
$ws.Range("N137").Value = -64200
# Row 139
$ws.Range("H139").Value = 1342.6666
$ws.Range("I139").Value = 1014
$ws.Range("K139").Value = 3042
$ws.Range("M139").Value = 2098
# Row 140
$ws.Range("H140").Value = 947.9
$ws.Range("I140").Value = 684.875
$ws.Range("K140").Value = 2054.625
$ws.Range("M140").Value = 3125.375

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 57
$ws.Range("H57").Value = 6266.6665
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
# Row 102
$ws.Range("H102").Value = 1549.0834
$ws.Range("I102").Value = 1618.9
$ws.Range("J102").Value = 1200
$ws.Range("K102").Value = 1618.9
$ws.Range("L102").Value = 1200
$ws.Range("M102").Value = 3.099999999999909
$ws.Range("N102").Value = -4444
# Row 107
$ws.Range("H107").Value = 345
$ws.Range("I107").Value = 290
$ws.Range("J107").Value = 400
$ws.Range("K107").Value = 290
$ws.Range("L107").Value = 400
$ws.Range("M107").Value = 1630
$ws.Range("N107").Value = -4240
# Row 122
$ws.Range("H122").Value = 3000
$ws.Range("I122").Value = 3000
$ws.Range("K122").Value = 9000
$ws.Range("M122").Value = -6550
# Row 132
$ws.Range("H132").Value = 2750
$ws.Range("I132").Value = 1000
$ws.Range("J132").Value = 4500
$ws.Range("K132").Value = 3000
$ws.Range("L132").Value = 13500
$ws.Range("M132").Value = -470
$ws.Range("N132").Value = -18560

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 2245.6
$ws.Range("I16").Value = 2309.5
$ws.Range("K16").Value = 2309.5
$ws.Range("M16").Value = -2139.5
# Row 40
$ws.Range("H40").Value = 5576
$ws.Range("I40").Value = 5576
$ws.Range("K40").Value = 5576
$ws.Range("M40").Value = -5440
# Row 41
$ws.Range("H41").Value = 20000
$ws.Range("J41").Value = 20000
$ws.Range("L41").Value = 20000
$ws.Range("N41").Value = -20876
# Row 42
$ws.Range("H42").Value = 40000000
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 40000000
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 40000000
$ws.Range("M42").ClearContents()
$ws.Range("N42").Value = -40001126
# Row 49
$ws.Range("H49").Value = 40000000
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 40000000
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 40000000
$ws.Range("M49").ClearContents()
$ws.Range("N49").Value = -40000294
# Row 55
$ws.Range("H55").Value = 1404.5
$ws.Range("J55").Value = 2002
$ws.Range("L55").Value = 2002
$ws.Range("N55").Value = -2348
# Row 8
$ws.Range("H8").Value = 1000
$ws.Range("J8").Value = 1000
$ws.Range("L8").Value = 1000
$ws.Range("N8").Value = -1280

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 76
$ws.Range("H76").Value = 58173
$ws.Range("J76").Value = 58173
$ws.Range("L76").Value = 58173
$ws.Range("N76").Value = -58803
# Row 79
$ws.Range("H79").Value = 58173
$ws.Range("J79").Value = 58173
$ws.Range("L79").Value = 58173
$ws.Range("N79").Value = -60357
# Row 122
$ws.Range("H122").Value = 2786.6667
$ws.Range("I122").Value = 2786.6667
$ws.Range("K122").Value = 8360.000100000001
$ws.Range("M122").Value = -5910.000100000001
# Row 136
$ws.Range("H136").Value = 1223.5
$ws.Range("I136").Value = 1223.5
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 3670.5
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -1120.5
$ws.Range("N136").ClearContents()

